# Update the build timestamp embedded in the "build_version" string from
# "February 03 2026 17.29.55 EST" to "February 03 2026 18.05.36 EST"
# across the About sheet (A2, A6) and the Boundaries and methane sources
# sheet (S2:S8).

$wb = $excel.ActiveWorkbook

$oldStamp = "February 03 2026 17.29.55 EST"
$newStamp = "February 03 2026 18.05.36 EST"

$aboutWs = $wb.Worksheets.Item("About")
$dataWs = $wb.Worksheets.Item("Boundaries and methane sources")

# --- About sheet ---
$a2 = $aboutWs.Range("A2")
$a2Text = $a2.Value()
$a2.Value = $a2Text.Replace($oldStamp, $newStamp)

$a6 = $aboutWs.Range("A6")
$a6Text = $a6.Value()
$a6.Value = $a6Text.Replace($oldStamp, $newStamp)

# --- Boundaries and methane sources sheet: column S, rows 2-8 ---
for ($row = 2; $row -le 8; $row++) {
    $cell = $dataWs.Cells.Item($row, 19)  # column S = 19
    $cellText = $cell.Value()
    $cell.Value = $cellText.Replace($oldStamp, $newStamp)
}
